$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing header row (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J for rows 2-15
$data = @{
    2  = @(1, 5)
    3  = @(1, 4)
    4  = @(1, 4)
    5  = @(1, 7)
    6  = @(3, 5)
    7  = @(8, 9)
    8  = @(1, 5)
    9  = @(1, 4)
    10 = @(1, 5)
    11 = @(1, 5)
    12 = @(1, 4)
    13 = @(4, 6)
    14 = @(6, 8)
    15 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
